$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "당기순이익(비지배)" (col J) and "자본총계(비지배)" (col O) figures
# for the 2014-2018 (rows 2-5) actuals; the source stopped reporting them.
foreach ($ref in @("J2","J3","J4","J5","O2","O3","O4","O5")) {
    $ws.Range($ref).ClearContents()
}

# Row 2
$ws.Range("D2").Value = 285576
$ws.Range("E2").Value = -2897
$ws.Range("F2").Value = -2897
$ws.Range("G2").Value = -3855
$ws.Range("H2").Value = -2878
$ws.Range("I2").Value = -2878
$ws.Range("K2").Value = 102557
$ws.Range("L2").Value = 53467
$ws.Range("M2").Value = 49090
$ws.Range("N2").Value = 49090
$ws.Range("P2").Value = 2915
$ws.Range("Q2").Value = 8792
$ws.Range("R2").Value = -9642
$ws.Range("S2").Value = 2119
$ws.Range("T2").Value = 9123
$ws.Range("U2").Value = -331
$ws.Range("V2").Value = 36376
$ws.Range("W2").Value = -1.01
$ws.Range("X2").Value = -1.01
$ws.Range("Y2").Value = -5.61
$ws.Range("Z2").Value = -2.6
$ws.Range("AA2").Value = 108.92
$ws.Range("AB2").Value = 1575.88
$ws.Range("AC2").Value = -2468
$ws.Range("AD2").Value = -19.61
$ws.Range("AE2").Value = 42166
$ws.Range("AF2").Value = 1.15
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 0.31
$ws.Range("AI2").Value = -6.1
$ws.Range("AJ2").Value = 112582792

# Row 3
$ws.Range("D3").Value = 178903
$ws.Range("E3").Value = 8176
$ws.Range("F3").Value = 8176
$ws.Range("G3").Value = 8127
$ws.Range("H3").Value = 6313
$ws.Range("I3").Value = 6313
$ws.Range("K3").Value = 107955
$ws.Range("L3").Value = 54056
$ws.Range("M3").Value = 53899
$ws.Range("N3").Value = 53899
$ws.Range("P3").Value = 2915
$ws.Range("Q3").Value = 24756
$ws.Range("R3").Value = -28019
$ws.Range("S3").Value = -1841
$ws.Range("T3").Value = 6544
$ws.Range("U3").Value = 18211
$ws.Range("V3").Value = 35952
$ws.Range("W3").Value = 4.57
$ws.Range("X3").Value = 3.53
$ws.Range("Y3").Value = 12.26
$ws.Range("Z3").Value = 6
$ws.Range("AA3").Value = 100.29
$ws.Range("AB3").Value = 1740.79
$ws.Range("AC3").Value = 5414
$ws.Range("AD3").Value = 14.67
$ws.Range("AE3").Value = 46297
$ws.Range("AF3").Value = 1.72
$ws.Range("AG3").Value = 2400
$ws.Range("AH3").Value = 3.02
$ws.Range("AI3").Value = 44.27
$ws.Range("AJ3").Value = 112582792

# Row 4
$ws.Range("D4").Value = 163218
$ws.Range("E4").Value = 16169
$ws.Range("F4").Value = 16169
$ws.Range("G4").Value = 15751
$ws.Range("H4").Value = 12054
$ws.Range("I4").Value = 12054
$ws.Range("K4").Value = 139590
$ws.Range("L4").Value = 75674
$ws.Range("M4").Value = 63916
$ws.Range("N4").Value = 63916
$ws.Range("P4").Value = 2915
$ws.Range("Q4").Value = 17220
$ws.Range("R4").Value = -19669
$ws.Range("S4").Value = 8111
$ws.Range("T4").Value = 10636
$ws.Range("U4").Value = 6584
$ws.Range("V4").Value = 47095
$ws.Range("W4").Value = 9.91
$ws.Range("X4").Value = 7.38
$ws.Range("Y4").Value = 20.46
$ws.Range("Z4").Value = 9.74
$ws.Range("AA4").Value = 118.4
$ws.Range("AB4").Value = 2084.62
$ws.Range("AC4").Value = 10337
$ws.Range("AD4").Value = 8.19
$ws.Range("AE4").Value = 54901
$ws.Range("AF4").Value = 1.54
$ws.Range("AG4").Value = 6200
$ws.Range("AH4").Value = 7.32
$ws.Range("AI4").Value = 59.89
$ws.Range("AJ4").Value = 112582792

# Row 5
$ws.Range("D5").Value = 208914
$ws.Range("E5").Value = 13733
$ws.Range("F5").Value = 13733
$ws.Range("G5").Value = 16449
$ws.Range("H5").Value = 12465
$ws.Range("I5").Value = 12465
$ws.Range("K5").Value = 150875
$ws.Range("L5").Value = 82448
$ws.Range("M5").Value = 68426
$ws.Range("N5").Value = 68426
$ws.Range("P5").Value = 2915
$ws.Range("Q5").Value = 11433
$ws.Range("R5").Value = -8322
$ws.Range("S5").Value = -5983
$ws.Range("T5").Value = 24141
$ws.Range("U5").Value = -12708
$ws.Range("V5").Value = 48442
$ws.Range("W5").Value = 6.57
$ws.Range("X5").Value = 5.97
$ws.Range("Y5").Value = 18.84
$ws.Range("Z5").Value = 8.58
$ws.Range("AA5").Value = 120.49
$ws.Range("AB5").Value = 2239.92
$ws.Range("AC5").Value = 10690
$ws.Range("AD5").Value = 10.94
$ws.Range("AE5").Value = 58775
$ws.Range("AF5").Value = 1.99
$ws.Range("AG5").Value = 5900
$ws.Range("AH5").Value = 5.04
$ws.Range("AI5").Value = 55.11
$ws.Range("AJ5").Value = 112582792

# Row 6
$ws.Range("D6").Value = 254633
$ws.Range("E6").Value = 6395
$ws.Range("F6").Value = 6395
$ws.Range("G6").Value = 3318
$ws.Range("H6").Value = 2580
$ws.Range("I6").Value = 2580
$ws.Range("K6").Value = 159550
$ws.Range("L6").Value = 94858
$ws.Range("M6").Value = 64692
$ws.Range("N6").Value = 64692
$ws.Range("P6").Value = 2915
$ws.Range("Q6").Value = -2882
$ws.Range("R6").Value = -4152
$ws.Range("S6").Value = 8878
$ws.Range("T6").Value = 20414
$ws.Range("U6").Value = -23296
$ws.Range("V6").Value = 64031
$ws.Range("W6").Value = 2.51
$ws.Range("X6").Value = 1.01
$ws.Range("Y6").Value = 3.88
$ws.Range("Z6").Value = 1.66
$ws.Range("AA6").Value = 146.63
$ws.Range("AB6").Value = 2111.25
$ws.Range("AC6").Value = 2213
$ws.Range("AD6").Value = 44.15
$ws.Range("AE6").Value = 55568
$ws.Range("AF6").Value = 1.76
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 0.77
$ws.Range("AI6").Value = 33.88
$ws.Range("AJ6").Value = 112582792

# Row 7
$ws.Range("D7").Value = 241711
$ws.Range("E7").Value = 6110
$ws.Range("G7").Value = 2741
$ws.Range("H7").Value = 2022
$ws.Range("I7").Value = 1955
$ws.Range("K7").Value = 165650
$ws.Range("L7").Value = 100083
$ws.Range("M7").Value = 65567
$ws.Range("N7").Value = 65499
$ws.Range("P7").Value = 2918
$ws.Range("Q7").Value = 8575
$ws.Range("R7").Value = -8767
$ws.Range("S7").Value = 4019
$ws.Range("T7").Value = 8087
$ws.Range("U7").Value = 410
$ws.Range("W7").Value = 2.53
$ws.Range("X7").Value = 0.84
$ws.Range("Y7").Value = 3
$ws.Range("Z7").Value = 1.24
$ws.Range("AA7").Value = 152.64
$ws.Range("AC7").Value = 1676
$ws.Range("AD7").Value = 48.26
$ws.Range("AE7").Value = 56260
$ws.Range("AF7").Value = 1.44
$ws.Range("AG7").Value = 719
$ws.Range("AH7").Value = 0.89
$ws.Range("AI7").Value = 41.41

# Row 8
$ws.Range("D8").Value = 252954
$ws.Range("E8").Value = 13032
$ws.Range("G8").Value = 12211
$ws.Range("H8").Value = 9164
$ws.Range("I8").Value = 8945
$ws.Range("K8").Value = 170887
$ws.Range("L8").Value = 98924
$ws.Range("M8").Value = 71963
$ws.Range("N8").Value = 71837
$ws.Range("P8").Value = 2918
$ws.Range("Q8").Value = 14255
$ws.Range("R8").Value = -8529
$ws.Range("S8").Value = -3274
$ws.Range("T8").Value = 8471
$ws.Range("U8").Value = 4200
$ws.Range("W8").Value = 5.15
$ws.Range("X8").Value = 3.62
$ws.Range("Y8").Value = 13.04
$ws.Range("Z8").Value = 5.46
$ws.Range("AA8").Value = 137.46
$ws.Range("AC8").Value = 7672
$ws.Range("AD8").Value = 9.93
$ws.Range("AE8").Value = 61704
$ws.Range("AF8").Value = 1.23
$ws.Range("AG8").Value = 3067
$ws.Range("AH8").Value = 4.03
$ws.Range("AI8").Value = 38.6

# Row 9
$ws.Range("D9").Value = 249526
$ws.Range("E9").Value = 13637
$ws.Range("G9").Value = 12863
$ws.Range("H9").Value = 9667
$ws.Range("I9").Value = 9443
$ws.Range("K9").Value = 175150
$ws.Range("L9").Value = 98109
$ws.Range("M9").Value = 77042
$ws.Range("N9").Value = 77018
$ws.Range("P9").Value = 2918
$ws.Range("Q9").Value = 16691
$ws.Range("R9").Value = -10086
$ws.Range("S9").Value = -5134
$ws.Range("T9").Value = 9994
$ws.Range("U9").Value = 4516
$ws.Range("W9").Value = 5.46
$ws.Range("X9").Value = 3.87
$ws.Range("Y9").Value = 12.69
$ws.Range("Z9").Value = 5.59
$ws.Range("AA9").Value = 127.34
$ws.Range("AC9").Value = 8098
$ws.Range("AD9").Value = 9.41
$ws.Range("AE9").Value = 66155
$ws.Range("AF9").Value = 1.15
$ws.Range("AG9").Value = 3512
$ws.Range("AH9").Value = 4.61
$ws.Range("AI9").Value = 41.88
